# Rename the forecast-year column headers from "f" (forecast) suffix to
# "p" (projection) suffix: "2026f" -> "2026p" and "2027f" -> "2027p".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value2 = "2026p"
$ws.Range("H1").Value2 = "2027p"
